$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.745.00'
$ws.Range("E2").Value = '  +2.31%  '
$ws.Range("D3").Value = '1.764.45'
$ws.Range("E3").Value = '  -0.91%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.54%  '
$ws.Range("D5").Value = '335.21'
$ws.Range("E5").Value = '  -0.23%  '
$ws.Range("D6").Value = '0.9985'
$ws.Range("E6").Value = '  -0.46%  '
$ws.Range("E7").Value = '  -0.34%  '
$ws.Range("D8").Value = '0.3424'
$ws.Range("E8").Value = '  +0.27%  '
$ws.Range("D9").Value = '46.65'
$ws.Range("E9").Value = '  -3.12%  '
$ws.Range("D10").Value = '1.133'
$ws.Range("E10").Value = '  -4.79%  '
$ws.Range("D11").Value = '0.07394'
$ws.Range("E11").Value = '  -0.94%  '
$ws.Range("D12").Value = '0.9975'
$ws.Range("E12").Value = '  -0.64%  '
$ws.Range("D13").Value = '22.30'
$ws.Range("E13").Value = '  +2.96%  '
$ws.Range("D14").Value = '6.337'
$ws.Range("E14").Value = '  -1.24%  '
$ws.Range("D15").Value = '1.764.57'
$ws.Range("E15").Value = '  -0.95%  '
$ws.Range("D16").Value = '7.054'
$ws.Range("E16").Value = '  -0.20%  '
$ws.Range("D17").Value = '0.00001073'
$ws.Range("E17").Value = '  -1.20%  '
$ws.Range("E18").Value = '  +0.16%  '
$ws.Range("D19").Value = '82.03'
$ws.Range("E19").Value = '  -1.59%  '
$ws.Range("D20").Value = '1.000'
$ws.Range("E20").Value = '  -0.29%  '
$ws.Range("D21").Value = '17.29'
$ws.Range("E21").Value = '  +0.11%  '
$ws.Range("D22").Value = '6.391'
$ws.Range("E22").Value = '  -2.61%  '
$ws.Range("D23").Value = '27.735.74'
$ws.Range("E23").Value = '  +2.27%  '
$ws.Range("D24").Value = '12.00'
$ws.Range("E24").Value = '  -1.92%  '
$ws.Range("D25").Value = '2.379'
$ws.Range("E25").Value = '  -0.27%  '
$ws.Range("D26").Value = '1.432'
$ws.Range("E26").Value = '  -2.44%  '
$ws.Range("D27").Value = '20.63'
$ws.Range("E27").Value = '  -2.47%  '
$ws.Range("D28").Value = '2.407'
$ws.Range("E28").Value = '  -4.30%  '
$ws.Range("D29").Value = '152.71'
$ws.Range("E29").Value = '  -0.94%  '
$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").Value = '134.08'
$ws.Range("E30").Value = '  +0.03%  '
$ws.Range("B31").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C31").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D31").Value = '1.961.31'
$ws.Range("E31").Value = '  -1.07%  '
$ws.Range("D32").Value = '6.088'
$ws.Range("E32").Value = '  +1.21%  '
$ws.Range("D33").Value = '3.959'
$ws.Range("E33").Value = '  -1.54%  '
$ws.Range("D34").Value = '0.08749'
$ws.Range("E34").Value = '  +0.83%  '
$ws.Range("D35").Value = '12.69'
$ws.Range("E35").Value = '  -3.51%  '
$ws.Range("E36").Value = '  +3.80%  '
$ws.Range("D37").Value = '0.6784'
$ws.Range("E37").Value = '  -0.88%  '
$ws.Range("D38").Value = '5.307'
$ws.Range("E38").Value = '  -1.50%  '
$ws.Range("D39").Value = '0.06281'
$ws.Range("E39").Value = '  -0.17%  '
$ws.Range("D40").Value = '0.2175'
$ws.Range("E40").Value = '  -0.08%  '
$ws.Range("D41").Value = '1.251'
$ws.Range("E41").Value = '  +1.45%  '
$ws.Range("D42").Value = '1.499'
$ws.Range("E42").Value = '  -8.09%  '
$ws.Range("D43").Value = '8.218'
$ws.Range("E43").Value = '  -5.88%  '
$ws.Range("D44").Value = '0.9984'
$ws.Range("E44").Value = '  -0.44%  '
$ws.Range("D45").Value = '14.02'
$ws.Range("E45").Value = '  -1.61%  '
$ws.Range("D46").Value = '0.6247'
$ws.Range("E46").Value = '  -2.75%  '
$ws.Range("E47").Value = '  -0.33%  '
$ws.Range("D48").Value = '131.55'
$ws.Range("E48").Value = '  +0.90%  '
$ws.Range("D49").Value = '2.068'
$ws.Range("E49").Value = '  -2.91%  '
$ws.Range("D50").Value = '0.07393'
$ws.Range("E50").Value = '  +3.99%  '
$ws.Range("D51").Value = '1.138'
$ws.Range("E51").Value = '  +2.15%  '
